$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 43.25
$ws.Range("E2").Value = 42.68999862670898
$ws.Range("F2").Value = 49.25
$ws.Range("G2").Value = 40.05099868774414
$ws.Range("H2").Value = 237933261
$ws.Range("I2").Value = "FRSH"

$ws.Range("D3").Value = 43.25
$ws.Range("E3").Value = 42.68999862670898
$ws.Range("F3").Value = 49.25
$ws.Range("G3").Value = 40.05099868774414
$ws.Range("H3").Value = 237933261
$ws.Range("I3").Value = "FRSH"

$ws.Range("D4").Value = 43.25
$ws.Range("E4").Value = 42.68999862670898
$ws.Range("F4").Value = 49.25
$ws.Range("G4").Value = 40.05099868774414
$ws.Range("H4").Value = 237933261
$ws.Range("I4").Value = "FRSH"

$ws.Range("D5").Value = 43.25
$ws.Range("E5").Value = 42.68999862670898
$ws.Range("F5").Value = 49.25
$ws.Range("G5").Value = 40.05099868774414
$ws.Range("H5").Value = 237933261
$ws.Range("I5").Value = "FRSH"

$ws.Range("D6").Value = 43.25
$ws.Range("E6").Value = 42.68999862670898
$ws.Range("F6").Value = 49.25
$ws.Range("G6").Value = 40.05099868774414
$ws.Range("H6").Value = 237933261
$ws.Range("I6").Value = "FRSH"

$ws.Range("D7").Value = 42.7400016784668
$ws.Range("E7").Value = 50.25
$ws.Range("F7").Value = 50.63000106811523
$ws.Range("G7").Value = 39
$ws.Range("H7").Value = 237933261
$ws.Range("I7").Value = "FRSH"

$ws.Range("D8").Value = 26.39999961853028
$ws.Range("E8").Value = 21.75
$ws.Range("F8").Value = 26.80999946594238
$ws.Range("G8").Value = 18.52000045776367
$ws.Range("H8").Value = 237933261
$ws.Range("I8").Value = "FRSH"

$ws.Range("D9").Value = 18.21999931335449
$ws.Range("E9").Value = 18.17000007629395
$ws.Range("F9").Value = 20.26000022888184
$ws.Range("G9").Value = 15.76000022888184
$ws.Range("H9").Value = 237933261
$ws.Range("I9").Value = "FRSH"

$ws.Range("D10").Value = 13.19999980926514
$ws.Range("E10").Value = 13.10999965667725
$ws.Range("F10").Value = 16.07999992370605
$ws.Range("G10").Value = 11.47000026702881
$ws.Range("H10").Value = 237933261
$ws.Range("I10").Value = "FRSH"

$ws.Range("D11").Value = 13.0600004196167
$ws.Range("E11").Value = 13.59000015258789
$ws.Range("F11").Value = 15.09000015258789
$ws.Range("G11").Value = 11.96000003814697
$ws.Range("H11").Value = 237933261
$ws.Range("I11").Value = "FRSH"

$ws.Range("D12").Value = 15.07999992370606
$ws.Range("E12").Value = 16.18000030517578
$ws.Range("F12").Value = 16.92499923706055
$ws.Range("G12").Value = 13.84000015258789
$ws.Range("H12").Value = 237933261
$ws.Range("I12").Value = "FRSH"

$ws.Range("D13").Value = 15
$ws.Range("E13").Value = 13.35999965667725
$ws.Range("F13").Value = 15.34500026702881
$ws.Range("G13").Value = 13.30000019073486
$ws.Range("H13").Value = 237933261
$ws.Range("I13").Value = "FRSH"

$ws.Range("D14").Value = 17.35000038146973
$ws.Range("E14").Value = 18.65999984741211
$ws.Range("F14").Value = 18.70999908447266
$ws.Range("G14").Value = 15.82499980926514
$ws.Range("H14").Value = 237933261
$ws.Range("I14").Value = "FRSH"

$ws.Range("D15").Value = 19.92000007629395
$ws.Range("E15").Value = 17.94000053405762
$ws.Range("F15").Value = 20.06999969482422
$ws.Range("G15").Value = 17.28499984741211
$ws.Range("H15").Value = 237933261
$ws.Range("I15").Value = "FRSH"

$ws.Range("D16").Value = 23.11400032043457
$ws.Range("E16").Value = 22.20000076293945
$ws.Range("F16").Value = 23.93000030517578
$ws.Range("G16").Value = 21.01499938964844
$ws.Range("H16").Value = 237933261
$ws.Range("I16").Value = "FRSH"

$ws.Range("D17").Value = 18.20000076293945
$ws.Range("E17").Value = 17.85000038146973
$ws.Range("F17").Value = 18.92000007629395
$ws.Range("G17").Value = 16.86000061035156
$ws.Range("H17").Value = 237933261
$ws.Range("I17").Value = "FRSH"

$ws.Range("D18").Value = 12.73999977111816
$ws.Range("E18").Value = 12.5
$ws.Range("F18").Value = 13.80000019073486
$ws.Range("G18").Value = 12.17500019073486
$ws.Range("H18").Value = 237933261
$ws.Range("I18").Value = "FRSH"

$ws.Range("D19").Value = 11.52000045776367
$ws.Range("E19").Value = 11.69999980926514
$ws.Range("F19").Value = 12
$ws.Range("G19").Value = 10.8100004196167
$ws.Range("H19").Value = 237933261
$ws.Range("I19").Value = "FRSH"

$ws.Range("D20").Value = 16.3799991607666
$ws.Range("E20").Value = 18.60000038146973
$ws.Range("F20").Value = 19.77000045776367
$ws.Range("G20").Value = 15.28499984741211
$ws.Range("H20").Value = 237933261
$ws.Range("I20").Value = "FRSH"

$ws.Range("D21").Value = 14.09000015258789
$ws.Range("E21").Value = 14.77000045776367
$ws.Range("F21").Value = 15.09000015258789
$ws.Range("G21").Value = 11.36999988555908
$ws.Range("H21").Value = 237933261
$ws.Range("I21").Value = "FRSH"

$ws.Range("D22").Value = 14.89000034332275
$ws.Range("E22").Value = 12.98999977111816
$ws.Range("F22").Value = 15.47000026702881
$ws.Range("G22").Value = 12.92500019073486
$ws.Range("H22").Value = 237933261
$ws.Range("I22").Value = "FRSH"

